# Update "想去人数" (F column) values on sheets "展览" and "全部类型"
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    3 = 26
    4 = 121
    5 = 108
    6 = 476
    9 = 602
    11 = 325
    12 = 26
    13 = 388
    17 = 55
    19 = 103
    20 = 1020
    21 = 1424
    22 = 310
    23 = 343
    25 = 83
    29 = 237
    30 = 265
    31 = 286
    34 = 102
    35 = 160
    36 = 594
    37 = 299
    38 = 3799
    39 = 2
    40 = 441
    41 = 218
    42 = 952
    43 = 51
    46 = 81
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}

Write-Output "Updated F column values on sheets: $($sheetNames -join ', ')"
